# batsmanset_smat.xlsx — add four new batsmen rows (JP Inglis, MP Breetzke,
# RD Rickelton, BJ Jacobs) pasted in from an external (HTML clipboard-style)
# source, which is why they carry their own font/fill formatting instead of
# the sheet's existing "body row" style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Cell values for the four new rows
# ---------------------------------------------------------------------

# Row 16 - JP Inglis
$ws.Range("A16").Value = "JP Inglis"
$ws.Range("B16").Value = "2017-2024"
$ws.Range("C16").Value = 78
$ws.Range("D16").Value = 72
$ws.Range("E16").Value = 11
$ws.Range("F16").Value = 1800
$ws.Range("G16").Value = 79
$ws.Range("H16").Value = 29.5
$ws.Range("I16").Value = 1284
$ws.Range("J16").Value = 140.18
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 14
$ws.Range("M16").Value = 3
$ws.Range("N16").Value = 170
$ws.Range("O16").Value = 56
$ws.Range("P16").Value = "BAT"
$ws.Range("Q16").Value = 5

# Row 17 - MP Breetzke
$ws.Range("A17").Value = "MP Breetzke"
$ws.Range("B17").Value = "2023-2025"
$ws.Range("C17").Value = 28
$ws.Range("D17").Value = 24
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 658
$ws.Range("G17").Value = 78
$ws.Range("H17").Value = 29.9
$ws.Range("I17").Value = 511
$ws.Range("J17").Value = 128.76
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 2
$ws.Range("M17").Value = 0
$ws.Range("N17").Value = 62
$ws.Range("O17").Value = 30
$ws.Range("P17").Value = "BAT"
$ws.Range("Q17").Value = 1

# Row 18 - RD Rickelton
$ws.Range("A18").Value = "RD Rickelton"
$ws.Range("B18").Value = "2023-2025"
$ws.Range("C18").Value = 25
$ws.Range("D18").Value = 25
$ws.Range("E18").Value = 2
$ws.Range("F18").Value = 1012
$ws.Range("G18").Value = 98
$ws.Range("H18").Value = 44
$ws.Range("I18").Value = 622
$ws.Range("J18").Value = 162.7
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 8
$ws.Range("M18").Value = 0
$ws.Range("N18").Value = 90
$ws.Range("O18").Value = 57
$ws.Range("P18").Value = "BAT"
$ws.Range("Q18").Value = 2

# Row 19 - BJ Jacobs (F19 feeds the E14 "=-F19" formula elsewhere on sheet)
$ws.Range("A19").Value = "BJ Jacobs"
$ws.Range("B19").Value = "2024-2025"
$ws.Range("C19").Value = 3
$ws.Range("D19").Value = 3
$ws.Range("E19").Value = 1
$ws.Range("F19").Value = 26
$ws.Range("G19").Value = 18
$ws.Range("H19").Value = 13
$ws.Range("I19").Value = 27
$ws.Range("J19").Value = 96.29
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = 0
$ws.Range("N19").Value = 2
$ws.Range("O19").Value = 1
$ws.Range("P19").Value = "BAT"
$ws.Range("Q19").Value = 5

# ---------------------------------------------------------------------
# 2. Formatting — rebuild the pasted-in look:
#    * A16:B18            -> dark grey (#222222) text, thin border, centered
#    * C16:O18,O19,Q16:Q19 -> same dark grey text + wrap, thin border, centered
#    * B19                -> medium grey (#48494A) text, thin border, centered (h only)
#    * C19:N19             -> medium grey (#48494A) text on white fill, thin border, centered
#    * P16:P19             -> thin border, centered (matches existing body style)
# ---------------------------------------------------------------------

$xlPasteFormats = -4122
$xlCenter = -4108

# --- style A: dark-grey font (#222222), border, center/center ---
$styleA = $ws.Range("A16")
$styleA.Font.Color = 2236962
$styleA.Borders.LineStyle = 1
$styleA.HorizontalAlignment = $xlCenter
$styleA.VerticalAlignment = $xlCenter

$styleA.Copy()
$ws.Range("B16").PasteSpecial($xlPasteFormats)
$styleA.Copy()
$ws.Range("A17:B18").PasteSpecial($xlPasteFormats)
$styleA.Copy()
$ws.Range("A19").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

# --- style B: dark-grey font (#222222) + wrap text, border, center/center ---
$styleB = $ws.Range("C16")
$styleB.Font.Color = 2236962
$styleB.Borders.LineStyle = 1
$styleB.WrapText = $true
$styleB.HorizontalAlignment = $xlCenter
$styleB.VerticalAlignment = $xlCenter

$styleB.Copy()
$ws.Range("D16:O16").PasteSpecial($xlPasteFormats)
$styleB.Copy()
$ws.Range("C17:O17").PasteSpecial($xlPasteFormats)
$styleB.Copy()
$ws.Range("C18:O18").PasteSpecial($xlPasteFormats)
$styleB.Copy()
$ws.Range("O19").PasteSpecial($xlPasteFormats)
$styleB.Copy()
$ws.Range("Q16:Q19").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

# --- style C: medium-grey font (#48494A), border, horizontal-center only ---
$styleC = $ws.Range("B19")
$styleC.Font.Color = 4868424
$styleC.Borders.LineStyle = 1
$styleC.HorizontalAlignment = $xlCenter
$excel.CutCopyMode = $false

# --- style D: medium-grey font (#48494A) on white fill, border, center/center ---
$styleD = $ws.Range("C19")
$styleD.Font.Color = 4868424
$styleD.Interior.Color = 16777215
$styleD.Borders.LineStyle = 1
$styleD.HorizontalAlignment = $xlCenter
$styleD.VerticalAlignment = $xlCenter

$styleD.Copy()
$ws.Range("D19:N19").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

# --- style for the "Type" column (P16:P19) — plain body style, centered ---
$pStyle = $ws.Range("P16")
$pStyle.Borders.LineStyle = 1
$pStyle.HorizontalAlignment = $xlCenter
$pStyle.VerticalAlignment = $xlCenter

$pStyle.Copy()
$ws.Range("P17:P19").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 3. Misc sheet/workbook bookkeeping to match the saved file
# ---------------------------------------------------------------------

# Selection left where the author's cursor ended up after entering data
$ws.Range("I23").Select()
